# edit.ps1 - apply "Half of presentation is ready." commit:
# Insert new slide-20 speech content (slides 21-25 + END OF PART 1 banner)
# between "Ok. Let's consider example." and the "Tell about ldd..." blurb,
# move the lastRenderedPageBreak marker there from the old position, and
# split the "k slaidu pro analiz vyvoda" run into per-word runs.

$d = $word.ActiveDocument

# --- 1. Insert the large new block of paragraphs (slides 21-25) -----------
# It replaces the single empty paragraph that used to sit right after
# "Ok. Let's consider example."
$anchor = $d.Content
$null = $anchor.Find.Execute("Ok. Let")
$anchor.Expand(4)
$anchor.Collapse(0)

$block1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">We have </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>struct</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Counter…</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:lastRenderedPageBreak/><w:t>Slide 21.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Two </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cpp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> files with static instances and a main function that simply prints the counter value.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>What you think will be the output?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>Slide 22.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">2? </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Do</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> someone has another opinion? Well, very reasonable.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>Slide 23.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">But now, let’s ship </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>A.o</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>B.o</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Counter.o</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as a static library</w:t></w:r><w:r><w:t xml:space="preserve">. Whoa. And the output now is 0. It is a little bit unexpected. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Add explanation here.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>Slide 24.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">If we will explicitly put </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>A.o</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the linkers input, than result will be one. This will not allow linker to throw unreferenced code away.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>Slide 25.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Also there is special option that forces linker to include </w:t></w:r><w:r><w:t>all object files from static library</w:t></w:r><w:r><w:t xml:space="preserve">. But be careful with it, as it can drastically increase the size of your binary. And do not forget to close it with </w:t></w:r><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Wl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>,-no-whole-archive</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>END OF PART 1.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$anchor.InsertXML($block1)

# --- 2. Move w:lastRenderedPageBreak onto the "Tell about ldd..." run -----
$anchor2 = $d.Content
$null = $anchor2.Find.Execute("Tell about")
$anchor2.Expand(4)

$tellXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Tell about </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>ldd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">In article on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>habr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$anchor2.InsertXML($tellXml)

# --- 3. ...and remove it from "Maybe find another name..." ----------------
$anchor3 = $d.Content
$null = $anchor3.Find.Execute("Maybe find another name")
$anchor3.Expand(4)

$maybeXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Maybe find another name for hereditary disease</w:t></w:r></w:p>
'@
$anchor3.InsertXML($maybeXml)

# --- 4. Split "к слайду про анализ вывода " into separate word runs -------
$anchor4 = $d.Content
$null = $anchor4.Find.Execute("Добавить")
$anchor4.Expand(4)

$para3Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Добавить</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t>к</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t>слайду</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">про анализ вывода </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>nm</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> для </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>static</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>member</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>mess</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>up</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>unix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>version</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> объяснение как раскрывается данная переменная во время загрузки </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>dll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$anchor4.InsertXML($para3Xml)
